$wb = $excel.ActiveWorkbook

# --- Deep water port sheet ---
$ws1 = $wb.Worksheets.Item("Deep water port")
$ws1.Range("D3").Formula = "=SUM(D4:D35)"
$ws1.Range("E3").Formula = "=SUM(E4:E35)"

# --- Shallow water port sheet ---
$ws2 = $wb.Worksheets.Item("Shallow water port")
$ws2.Range("D3").Formula = "=SUM(D4:D21)"
$ws2.Range("E3").Formula = "=SUM(E4:E21)"

# Corrected ship / BR listing for rows 4-21
$ws2.Cells.Item(4, 2).Value = "Hercules"
$ws2.Cells.Item(4, 3).Value = 100

$ws2.Cells.Item(5, 2).Value = "Pandora"
$ws2.Cells.Item(5, 3).Value = 100

$ws2.Cells.Item(6, 2).Value = "Mercury"
$ws2.Cells.Item(6, 3).Value = 80

$ws2.Cells.Item(7, 2).Value = "Mortar Brig"
$ws2.Cells.Item(7, 3).Value = 80

$ws2.Cells.Item(8, 2).Value = "NavyBrig"
$ws2.Cells.Item(8, 3).Value = 80

$ws2.Cells.Item(9, 2).Value = "Niagara"
$ws2.Cells.Item(9, 3).Value = 80

$ws2.Cells.Item(10, 2).Value = "Prince de Neufchatel"
$ws2.Cells.Item(10, 3).Value = 80

$ws2.Cells.Item(11, 2).Value = "Rattlesnake"
$ws2.Cells.Item(11, 3).Value = 80

$ws2.Cells.Item(12, 2).Value = "Rattlesnake Heavy"
$ws2.Cells.Item(12, 3).Value = 80

$ws2.Cells.Item(13, 2).Value = "Snow"
$ws2.Cells.Item(13, 3).Value = 80

$ws2.Cells.Item(14, 2).Value = "Brig"
$ws2.Cells.Item(14, 3).Value = 70

$ws2.Cells.Item(15, 2).Value = "Pickle"
$ws2.Cells.Item(15, 3).Value = 55

$ws2.Cells.Item(16, 2).Value = "Cutter"
$ws2.Cells.Item(16, 3).Value = 50

$ws2.Cells.Item(17, 2).Value = "GunBoat"
$ws2.Cells.Item(17, 3).Value = 50

$ws2.Cells.Item(18, 2).Value = "Lynx"
$ws2.Cells.Item(18, 3).Value = 50

$ws2.Cells.Item(19, 2).Value = "Privateer"
$ws2.Cells.Item(19, 3).Value = 50

$ws2.Cells.Item(20, 2).Value = "Yacht"
$ws2.Cells.Item(20, 3).Value = 50

$ws2.Cells.Item(21, 2).Value = "Yacht Silver"
$ws2.Cells.Item(21, 3).Value = 50
